$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("product backlog")
$wsSprint  = $wb.Worksheets.Item("sprint backlog")

# ---------------------------------------------------------------------
# "product backlog" sheet: rows 6 and 10 become "completed" items, so
# mirror the grey/yellow highlighted formatting used by the other
# completed/non-functional rows (row 4 is the template for this look).
# ---------------------------------------------------------------------
$wsProduct.Range("A4:I4").Copy() | Out-Null
$wsProduct.Range("A6:I6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsProduct.Range("I6").Value = "(3) Completed"

$wsProduct.Range("A7:I7").Copy() | Out-Null
$wsProduct.Range("A10:I10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (row 7 already carries ht=30, s=5/customFormat look)
$wsProduct.Range("I10").Value = "(3) Completed"

# ---------------------------------------------------------------------
# "sprint backlog" sheet: row 2 gets replaced with a different backlog
# item (id 8 / "to properly insert the organizational chart in the
# database"), and loses its Theme/Notes values.
# ---------------------------------------------------------------------
$wsSprint.Range("B2").Clear() | Out-Null
$wsSprint.Range("G2").Clear() | Out-Null
$wsSprint.Range("A2").Value = 8
$wsSprint.Range("E2").Value = "to properly insert the organizational chart in the database"

# Row grows to the taller "30" auto-height used elsewhere for wrapped,
# single-line-note rows (mirrors product-backlog row 8).
$wsSprint.Range("A2:J2").RowHeight = 30

# ---------------------------------------------------------------------
# Selections / active sheet & cells
# ---------------------------------------------------------------------
$wsProduct.Range("A10:XFD10").Select() | Out-Null
$wsSprint.Activate() | Out-Null
$wsSprint.Range("G5").Select() | Out-Null
